$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (pushes existing rows 20-111 down to 21-112)
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new daily record (day 19, 08/2025)
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 15547.52
$ws.Cells.Item(20, 3).Value = 8
$ws.Cells.Item(20, 4).Value = 2025
$ws.Cells.Item(20, 5).Value = "08/2025"
